# Applies the BOQ (Bill of Quantities) row updates described in the commit diff.
# Columns B, C, F, H hold genuine numbers; columns A, D, E, G, I hold text
# (the sheet keeps "S. No."/rate-total columns as text-formatted numbers, flagged
# by Excel as "Number Stored as Text"). A leading `'` reproduces that same
# quote-prefixed-text behaviour for any new value that looks numeric or is blank,
# exactly like a user typing `'3` into a text-formatted cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: P. point / Short point -> Medium point (up to 6 mtr.)
$ws.Range("C8").Value = 82
$ws.Range("D8").Value = '''3'
$ws.Range("E8").Value = 'Medium point (up to 6 mtr.)'
$ws.Range("F8").Value = 472
$ws.Range("G8").Value = '''38704.00'

# Row 9: (blank) -> P. point / Rewiring text -> Long point (up to 10 mtr.)
$ws.Range("A9").Value = 'P. point'
$ws.Range("C9").Value = 64
$ws.Range("D9").Value = '''4'
$ws.Range("E9").Value = 'Long point  (up to 10 mtr.)'
$ws.Range("F9").Value = 662
$ws.Range("G9").Value = '''42368.00'

# Row 10: Each / 3/5 pin 6A socket -> 3/6 pin 16A socket
$ws.Range("C10").Value = 56
$ws.Range("D10").Value = '''6.0'
$ws.Range("E10").Value = 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F10").Value = 78
$ws.Range("G10").Value = '''4368.00'

# Row 11: (blank) -> R. mtr. / ceiling fan text -> 25 mm
$ws.Range("A11").Value = 'R. mtr.'
$ws.Range("C11").Value = 27
$ws.Range("D11").Value = '''17'
$ws.Range("E11").Value = '25 mm'
$ws.Range("F11").Value = 56
$ws.Range("G11").Value = '''1512.00'

# Row 12: Each -> Set / 1200mm sweep fan -> Plate earthing text
$ws.Range("A12").Value = 'Set'
$ws.Range("C12").Value = 48
$ws.Range("D12").Value = '''13.0'
$ws.Range("E12").Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F12").Value = 5733
$ws.Range("G12").Value = '''275184.00'

# Row 13: Each -> (blank) / LED batten (1170mm) -> LED batten (IP20 SMD) text
$ws.Range("A13").Value = ''''
$ws.Range("C13").Value = 96
$ws.Range("D13").Value = '''16.0'
$ws.Range("E13").Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = '''0.00'

# Row 14: 6A to 32A rating -> 50/63A rating
$ws.Range("C14").Value = 37
$ws.Range("D14").Value = '''32'
$ws.Range("E14").Value = ' 50/63 A rating'
$ws.Range("F14").Value = 900
$ws.Range("G14").Value = '''33300.00'

# Row 15: (blank) -> Each / Metal door text -> 8 Way (8+2)
$ws.Range("A15").Value = 'Each'
$ws.Range("C15").Value = 13
$ws.Range("D15").Value = '''35'
$ws.Range("E15").Value = '8 Way (8+2)'
$ws.Range("F15").Value = 2184
$ws.Range("G15").Value = '''28392.00'

# Row 16: % -> (blank) / Add Tender Premium -> Total
$ws.Range("A16").Value = ''''
$ws.Range("C16").Value = 38
$ws.Range("D16").Value = '''36'
$ws.Range("E16").Value = 'Total'

# Row 17: Qty executed upto date only
$ws.Range("C17").Value = 35

# Row 19: Grand Total Rs. amounts
$ws.Range("G19").Value = '''423828.00'
$ws.Range("H19").Value = '''423828.00'

# Row 21: NET PAYABLE AMOUNT Rs. amounts
$ws.Range("G21").Value = '''423828.00'
$ws.Range("H21").Value = '''423828.00'
